$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 234.11111
$ws.Cells.Item(33, 9).Value = 241.6
$ws.Cells.Item(33, 11).Value = 241.6
$ws.Cells.Item(33, 13).Value = -12.59999999999999
$ws.Cells.Item(40, 8).Value = 5251
$ws.Cells.Item(40, 9).Value = 6501
$ws.Cells.Item(40, 11).Value = 6501
$ws.Cells.Item(40, 13).Value = -6326
$ws.Cells.Item(64, 8).Value = 5400.6665
$ws.Cells.Item(64, 9).Value = 5025.25
$ws.Cells.Item(64, 11).Value = 5025.25
$ws.Cells.Item(64, 13).Value = -4777.25
$ws.Cells.Item(67, 8).Value = 5400.6665
$ws.Cells.Item(67, 9).Value = 5025.25
$ws.Cells.Item(67, 11).Value = 5025.25
$ws.Cells.Item(67, 13).Value = -4167.25
$ws.Cells.Item(70, 8).Value = 1295822.2
$ws.Cells.Item(70, 10).Value = 1457700.8
$ws.Cells.Item(70, 12).Value = 4373102.4
$ws.Cells.Item(70, 14).Value = -4373642.4
$ws.Cells.Item(73, 8).Value = 1295822.2
$ws.Cells.Item(73, 10).Value = 1457700.8
$ws.Cells.Item(73, 12).Value = 4373102.4
$ws.Cells.Item(73, 14).Value = -4374974.4
$ws.Cells.Item(76, 8).Value = 20014798
$ws.Cells.Item(76, 10).Value = 6999.5
$ws.Cells.Item(76, 12).Value = 6999.5
$ws.Cells.Item(76, 14).Value = -7629.5
$ws.Cells.Item(79, 8).Value = 20014798
$ws.Cells.Item(79, 10).Value = 6999.5
$ws.Cells.Item(79, 12).Value = 6999.5
$ws.Cells.Item(79, 14).Value = -9183.5
$ws.Cells.Item(80, 8).Value = 1661.6111
$ws.Cells.Item(80, 9).Value = 765.2
$ws.Cells.Item(80, 10).Value = 2006.3846
$ws.Cells.Item(80, 11).Value = 2295.6
$ws.Cells.Item(80, 12).Value = 6019.1538
$ws.Cells.Item(80, 13).Value = -1297.6
$ws.Cells.Item(80, 14).Value = -8015.1538
$ws.Cells.Item(83, 8).Value = 1661.6111
$ws.Cells.Item(83, 9).Value = 765.2
$ws.Cells.Item(83, 10).Value = 2006.3846
$ws.Cells.Item(83, 11).Value = 6886.8
$ws.Cells.Item(83, 12).Value = 18057.4614
$ws.Cells.Item(83, 13).Value = -1894.8
$ws.Cells.Item(83, 14).Value = -28041.4614
$ws.Cells.Item(100, 8).Value = 2440.85
$ws.Cells.Item(100, 9).Value = 1935
$ws.Cells.Item(100, 11).Value = 1935
$ws.Cells.Item(100, 13).Value = -1394
$ws.Cells.Item(129, 8).Value = 1891.9375
$ws.Cells.Item(129, 9).Value = 874
$ws.Cells.Item(129, 10).Value = 2683.6667
$ws.Cells.Item(129, 11).Value = 2622
$ws.Cells.Item(129, 12).Value = 8051.000100000001
$ws.Cells.Item(129, 13).Value = 2378
$ws.Cells.Item(129, 14).Value = -18051.0001
$ws.Cells.Item(132, 8).Value = 1456.2106
$ws.Cells.Item(132, 9).Value = 1474.6471
$ws.Cells.Item(132, 10).Value = 1299.5
$ws.Cells.Item(132, 11).Value = 4423.9413
$ws.Cells.Item(132, 12).Value = 3898.5
$ws.Cells.Item(132, 13).Value = -1893.9413
$ws.Cells.Item(132, 14).Value = -8958.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 23492.473
$ws.Cells.Item(32, 10).Value = 20998.5
$ws.Cells.Item(32, 12).Value = 20998.5
$ws.Cells.Item(32, 14).Value = -21572.5
$ws.Cells.Item(45, 8).Value = 1599.3125
$ws.Cells.Item(45, 9).Value = 1471.8334
$ws.Cells.Item(45, 10).Value = 1981.75
$ws.Cells.Item(45, 11).Value = 1471.8334
$ws.Cells.Item(45, 12).Value = 1981.75
$ws.Cells.Item(45, 13).Value = -1094.8334
$ws.Cells.Item(45, 14).Value = -2735.75
$ws.Cells.Item(63, 8).Value = 7179.484
$ws.Cells.Item(63, 9).Value = 2255.4
$ws.Cells.Item(63, 10).Value = 9524.286
$ws.Cells.Item(63, 11).Value = 2255.4
$ws.Cells.Item(63, 12).Value = 9524.286
$ws.Cells.Item(63, 13).Value = -1569.4
$ws.Cells.Item(63, 14).Value = -10896.286
$ws.Cells.Item(66, 8).Value = 7179.484
$ws.Cells.Item(66, 9).Value = 2255.4
$ws.Cells.Item(66, 10).Value = 9524.286
$ws.Cells.Item(66, 11).Value = 11277
$ws.Cells.Item(66, 12).Value = 47621.43
$ws.Cells.Item(66, 13).Value = -7845
$ws.Cells.Item(66, 14).Value = -54485.43
$ws.Cells.Item(74, 8).Value = 304797.84
$ws.Cells.Item(74, 10).Value = 3394.1
$ws.Cells.Item(74, 12).Value = 3394.1
$ws.Cells.Item(74, 14).Value = -5142.1
$ws.Cells.Item(77, 8).Value = 304797.84
$ws.Cells.Item(77, 10).Value = 3394.1
$ws.Cells.Item(77, 12).Value = 16970.5
$ws.Cells.Item(77, 14).Value = -25706.5
$ws.Cells.Item(88, 8).Value = 1219.25
$ws.Cells.Item(88, 9).Value = 1374.6666
$ws.Cells.Item(88, 10).Value = 1167.4445
$ws.Cells.Item(88, 11).Value = 1374.6666
$ws.Cells.Item(88, 12).Value = 1167.4445
$ws.Cells.Item(88, 13).Value = -968.6666
$ws.Cells.Item(88, 14).Value = -1979.4445
$ws.Cells.Item(91, 8).Value = 1219.25
$ws.Cells.Item(91, 9).Value = 1374.6666
$ws.Cells.Item(91, 10).Value = 1167.4445
$ws.Cells.Item(91, 11).Value = 1374.6666
$ws.Cells.Item(91, 12).Value = 1167.4445
$ws.Cells.Item(91, 13).Value = 29.33339999999998
$ws.Cells.Item(91, 14).Value = -3975.4445
$ws.Cells.Item(122, 8).Value = 1998.8572
$ws.Cells.Item(122, 9).Value = 1996.5
$ws.Cells.Item(122, 11).Value = 5989.5
$ws.Cells.Item(122, 13).Value = -3539.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 1914.2
$ws.Cells.Item(26, 9).Value = 1914.2
$ws.Cells.Item(26, 11).Value = 1914.2
$ws.Cells.Item(26, 13).Value = -1622.2
$ws.Cells.Item(53, 8).Value = 84000
$ws.Cells.Item(53, 10).Value = 84000
$ws.Cells.Item(53, 12).Value = 84000
$ws.Cells.Item(53, 14).Value = -85148
$ws.Cells.Item(107, 8).Value = 1354.15
$ws.Cells.Item(107, 9).Value = 852.7692
$ws.Cells.Item(107, 10).Value = 2285.2856
$ws.Cells.Item(107, 11).Value = 852.7692
$ws.Cells.Item(107, 12).Value = 2285.2856
$ws.Cells.Item(107, 13).Value = 1067.2308
$ws.Cells.Item(107, 14).Value = -6125.2856
$ws.Cells.Item(134, 8).Value = 6314.143
$ws.Cells.Item(134, 10).Value = 9506.1
$ws.Cells.Item(134, 12).Value = 28518.3
$ws.Cells.Item(134, 14).Value = -33588.3
$ws.Cells.Item(137, 8).Value = 119758
$ws.Cells.Item(137, 10).Value = 119758
$ws.Cells.Item(137, 12).Value = 119758
$ws.Cells.Item(137, 14).Value = -129958

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 9128.817999999999
$ws.Cells.Item(62, 9).Value = 6788.143
$ws.Cells.Item(62, 10).Value = 13225
$ws.Cells.Item(62, 11).Value = 6788.143
$ws.Cells.Item(62, 12).Value = 13225
$ws.Cells.Item(62, 13).Value = -6164.143
$ws.Cells.Item(62, 14).Value = -14473
$ws.Cells.Item(65, 8).Value = 9128.817999999999
$ws.Cells.Item(65, 9).Value = 6788.143
$ws.Cells.Item(65, 10).Value = 13225
$ws.Cells.Item(65, 11).Value = 33940.715
$ws.Cells.Item(65, 12).Value = 66125
$ws.Cells.Item(65, 13).Value = -30820.715
$ws.Cells.Item(65, 14).Value = -72365
$ws.Cells.Item(130, 8).Value = 54428.57
$ws.Cells.Item(130, 10).Value = 54428.57
$ws.Cells.Item(130, 12).Value = 54428.57
$ws.Cells.Item(130, 14).Value = -64468.57

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(116, 8).Value = 4257
$ws.Cells.Item(116, 9).Value = 2029
$ws.Cells.Item(116, 10).Value = 4999.6665
$ws.Cells.Item(116, 11).Value = 6087
$ws.Cells.Item(116, 12).Value = 14998.9995
$ws.Cells.Item(116, 13).Value = -2645
$ws.Cells.Item(116, 14).Value = -21882.9995
$ws.Cells.Item(117, 8).Value = 1618.8572
$ws.Cells.Item(117, 10).Value = 2046.6
$ws.Cells.Item(117, 12).Value = 6139.799999999999
$ws.Cells.Item(117, 14).Value = -13023.8
$ws.Cells.Item(129, 8).Value = 26316812
$ws.Cells.Item(129, 10).Value = 55557196
$ws.Cells.Item(129, 12).Value = 166671588
$ws.Cells.Item(129, 14).Value = -166681588
$ws.Cells.Item(140, 8).Value = 1416.25
$ws.Cells.Item(140, 9).Value = 1331.6666
$ws.Cells.Item(140, 10).Value = 1444.4445
$ws.Cells.Item(140, 11).Value = 3994.9998
$ws.Cells.Item(140, 12).Value = 4333.333500000001
$ws.Cells.Item(140, 13).Value = 1185.0002
$ws.Cells.Item(140, 14).Value = -14693.3335

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 16261
$ws.Cells.Item(39, 10).Value = 16261
$ws.Cells.Item(39, 12).Value = 16261
$ws.Cells.Item(39, 14).Value = -17325
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()
$ws.Cells.Item(97, 8).Value = 970.3570999999999
$ws.Cells.Item(97, 10).Value = 1184.0714
$ws.Cells.Item(97, 12).Value = 1184.0714
$ws.Cells.Item(97, 14).Value = -2176.0714

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 27781160
$ws.Cells.Item(100, 9).Value = 125001850
$ws.Cells.Item(100, 11).Value = 125001850
$ws.Cells.Item(100, 13).Value = -125001309
$ws.Cells.Item(136, 8).Value = 5563471
$ws.Cells.Item(136, 9).Value = 7414317
$ws.Cells.Item(136, 11).Value = 22242951
$ws.Cells.Item(136, 13).Value = -22240401

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(101, 8).Value = 28799
$ws.Cells.Item(101, 10).Value = 28799
$ws.Cells.Item(101, 12).Value = 28799
$ws.Cells.Item(101, 14).Value = -35289
$ws.Cells.Item(107, 8).Value = 1683.1428
$ws.Cells.Item(107, 9).Value = 1235.5883
$ws.Cells.Item(107, 11).Value = 3706.7649
$ws.Cells.Item(107, 13).Value = -1786.7649
